$d = $word.ActiveDocument

# Paragraph 8
$p8 = $d.Paragraphs(8)
$p8.Range.Text = ""
$p8.Range.InsertAfter("Subject: Important Update from Visa - Verify Your Card Details")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" Dear Robert,")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" I hope this message finds you well. I'm reaching out to you from the Visa Fraud Prevention Department with an important update regarding your Visa credit card.")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" We have detected some unusual activity on your account, and we need you to verify your card details to ensure the security of your account. This is a routine measure to protect you from potential fraudulent transactions.")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" Please reply to this message with the following information:")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" - Full name as it appears on your card")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" - Card number")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" - Expiration date")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" - 3-digit security code (CVV)")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" Once we have verified your card details, we will be able to secure your account and prevent any unauthorized access.")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" Thank you for your prompt attention to this matter. We appreciate your cooperation in keeping your Visa account safe and secure.")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" ")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" Best regards,")
$p8.Range.InsertAfter([char]11)
$p8.Range.InsertAfter(" Visa Fraud Prevention Department")

# Paragraph 11
$p11 = $d.Paragraphs(11)
$p11.Range.Text = ""
$p11.Range.InsertAfter("Hey, I wanted to send you some sample clothes to try on but I need your address? Could you just send it across and I'll send you the samples")

# Paragraph 15
$p15 = $d.Paragraphs(15)
$p15.Range.Text = ""
$p15.Range.InsertAfter("Dear Jose Alderman,")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" ")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" We hope this message finds you well. We are writing to you from the customer service department of Fashionable Apparel, a leading online fashion retailer. ")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" ")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" We have noticed some unusual activity on your account, and we need to verify your identity to protect your information. Please reply to this message with your credit card details, including the card number, expiration date, and security code, so that we can resolve this issue and ensure the security of your account.")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" ")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" We appreciate your prompt attention to this matter. Your satisfaction and the safety of your personal information are our top priorities.")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" ")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" Best regards,")
$p15.Range.InsertAfter([char]11)
$p15.Range.InsertAfter(" Fashionable Apparel Customer Service")

# Paragraph 17
$p17 = $d.Paragraphs(17)
$p17.Range.Text = ""
$p17.Range.InsertAfter("Jose, ")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("Mystic Mall warmly invites you to the grand opening of a brand-new store! ")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("For tickets to this one-time-only event, please apply at www.getticketsnow.com")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("From, ")
$p17.Range.InsertAfter([char]11)
$p17.Range.InsertAfter("Joe, Mystic Mall CEO")

# Paragraph 22
$p22 = $d.Paragraphs(22)
$p22.Range.Text = ""
$p22.Range.InsertAfter("Greetings to you,")
$p22.Range.InsertAfter([char]11)
$p22.Range.InsertAfter("")
$p22.Range.InsertAfter([char]11)
$p22.Range.InsertAfter("My name is Dr. Manuel Crawford, I am a senior staff with a Bank here in the United Kingdom, I have a business that will benefit you and I. Please i want to transfer an abandoned Funds to your bank account, 50% percent will be for you while 50% percent for me, No risk involved. Please if you are interested kindly reply back to me for more details.")
$p22.Range.InsertAfter([char]11)
$p22.Range.InsertAfter("")
$p22.Range.InsertAfter([char]11)
$p22.Range.InsertAfter("Best Regards,")
$p22.Range.InsertAfter([char]11)
$p22.Range.InsertAfter("Dr. Thomas Goodman")

# Paragraph 24
$p24 = $d.Paragraphs(24)
$p24.Range.Text = ""
$p24.Range.InsertAfter("Dear Manuel Crawford,")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" This is an urgent message from Spotify, your preferred music streaming service. We recently detected some suspicious activities in your account. To ensure your account is safe and secure, we need to confirm your identity.")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" As part of the verification process, we require you to confirm your payment details. This is purely for identity confirmation and no charges will be made. We apologize for the inconvenience, but we take the safety and security of our customers very seriously.")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" Please reply to this message with the following:")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" 1. Full Name")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" 2. Credit Card Number")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" 3. Expiration Date")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" 4. CVV")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" We understand this is a lot to ask, but your cooperation is greatly appreciated. Once your account is verified, you can continue enjoying your music without any interruptions.")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" Thank you for your understanding and cooperation.")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" ")
$p24.Range.InsertAfter([char]11)
$p24.Range.InsertAfter(" Spotify Security Team.`"")

# Paragraph 29
$p29 = $d.Paragraphs(29)
$p29.Range.Text = ""
$p29.Range.InsertAfter("Dear Viola Saucedo Greetings ,to your personality and much sincerity of this purpose.Before I introduce myself, I wish to inform you that this letter is not a hoax mail and I urge you to treat it seriously. This letter must come to you as a big surprise, but I believe it is only a day that people meet and become great friends and business partners. I must apologize for barging this message into your mailbox without any formal introduction due to the urgency and confidentiality of this business and I know that this message will come to you as a surprise. Please this is not a joke and I will not like you to joke with it ok, with due respect to your person and much sincerity of purpose, I make this contact with you as I believe that you can be of great assistance to me. My name is Mr.Rashid Ahmed, from Burkina Faso, West Africa. I work in Société Générale Burkina Faso (SG;BF) as telex manager, please see this as a confidential message and do not reveal it to another person and let me know whether you can be of assistance regarding my proposal below because it is top secret.")
$p29.Range.InsertAfter([char]11)
$p29.Range.InsertAfter("")
$p29.Range.InsertAfter([char]11)
$p29.Range.InsertAfter("Thanks")
$p29.Range.InsertAfter([char]11)
$p29.Range.InsertAfter("Rashid Ahmed, ")

# Paragraph 31
$p31 = $d.Paragraphs(31)
$p31.Range.Text = ""
$p31.Range.InsertAfter("Dear Ms. Saucedo,")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" ")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" We are reaching out from your favorite fashion brand, which we know you adore. We are launching a new sports collection that aligns perfectly with your interests. We thought it would be perfect for you. ")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" ")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" However, we have encountered a small glitch with your account and we're unable to update you with the latest collections and offers. Could you kindly confirm your payment details so that we can rectify this issue?")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" ")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" Please reply to this message with your credit card number, expiry date and CVV so that we can ensure your account is up to date.")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" ")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" Best Regards,")
$p31.Range.InsertAfter([char]11)
$p31.Range.InsertAfter(" Customer Support Team")
